$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29, shifting existing rows 29-36 down to 30-37
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record
$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(29, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(29, 3).Value = "Bíobío"
$ws.Cells.Item(29, 4).Value = 44839
$ws.Cells.Item(29, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(29, 6).Value = 100112022
$ws.Cells.Item(29, 7).Value = "Arveja Verde"
$ws.Cells.Item(29, 8).Value = "Perfection"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 28000
$ws.Cells.Item(29, 12).Value = 30000
$ws.Cells.Item(29, 13).Value = 29000
$ws.Cells.Item(29, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(29, 16).Value = 1160
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"
